# Update BOC USD rates (auto)
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "All Published Values" - append new row 25 ---
$ws1 = $wb.Worksheets.Item("All Published Values")

$row25 = $ws1.Range("A25:J25")

# Force text storage (source values are number/date-looking strings that must
# stay text, matching the original inlineStr cells) then drop back to the
# default "Normal" style so the cells don't carry a lingering text format.
$row25.NumberFormat = "@"

$ws1.Range("A25").Value = "2026-01-04"
$ws1.Range("B25").Value = "2026-01-04 05:30:00"
$ws1.Range("C25").Value = "697.85"
$ws1.Range("D25").Value = "697.85"
$ws1.Range("E25").Value = "700.79"
$ws1.Range("F25").Value = "700.79"
$ws1.Range("G25").Value = "702.88"
$ws1.Range("H25").Value = "2026/01/04 05:30:00"
$ws1.Range("I25").Value = "2026-01-03 21:37:58"
$ws1.Range("J25").Value = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"

$row25.Style = "Normal"

# Extend the autofilter range to cover the new row
$ws1.AutoFilterMode = $false
$ws1.Range("A1:J25").AutoFilter() | Out-Null

# Refresh the hidden _FilterDatabase defined name so it also covers the new row
$wb.Names.Item("All Published Values!_FilterDatabase").RefersTo = "='All Published Values'!`$A`$1:`$J`$25"

# --- Sheet 2: "Daily Summary" - update publishes count for 2026-01-04 ---
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Range("B6").Value = 2
